$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 287-298 (only the cells that actually changed) ---

# Row 287 (Murcott / Especial) -> new date block 44448
$ws.Range("D287").Value = 44448
$ws.Range("M287").Value = 400
$ws.Range("N287").Value = 5000
$ws.Range("O287").Value = 5500
$ws.Range("P287").Value = 5250
$ws.Range("S287").Value = 525

# Row 288 (Murcott / Primera) -> new date block 44448
$ws.Range("D288").Value = 44448
$ws.Range("M288").Value = 600
$ws.Range("N288").Value = 4000
$ws.Range("O288").Value = 4500
$ws.Range("P288").Value = 4250
$ws.Range("S288").Value = 425

# Row 289 (Murcott / Segunda) -> new date block 44448
$ws.Range("D289").Value = 44448
$ws.Range("M289").Value = 500
$ws.Range("N289").Value = 3000
$ws.Range("O289").Value = 3500
$ws.Range("P289").Value = 3250
$ws.Range("S289").Value = 325

# Row 290 (Murcott / Tercera) -> new date block 44448
$ws.Range("D290").Value = 44448
$ws.Range("M290").Value = 300
$ws.Range("N290").Value = 2000
$ws.Range("O290").Value = 2500
$ws.Range("P290").Value = 2250
$ws.Range("S290").Value = 225

# Row 291 (was Clemenuless/44399 -> becomes Murcott/44167, Especial)
$ws.Range("D291").Value = 44167
$ws.Range("K291").Value = "Murcott"
$ws.Range("M291").Value = 360
$ws.Range("N291").Value = 8000
$ws.Range("O291").Value = 8500
$ws.Range("P291").Value = 8250
$ws.Range("S291").Value = 825

# Row 292 (was Clemenuless/44399 -> becomes Murcott/44167, Primera)
$ws.Range("D292").Value = 44167
$ws.Range("K292").Value = "Murcott"
$ws.Range("M292").Value = 280
$ws.Range("N292").Value = 7000
$ws.Range("O292").Value = 7500
$ws.Range("P292").Value = 7250
$ws.Range("S292").Value = 725

# Row 293 (was Clemenuless/44399 -> becomes Murcott/44167, Segunda)
$ws.Range("D293").Value = 44167
$ws.Range("K293").Value = "Murcott"
$ws.Range("M293").Value = 240
$ws.Range("N293").Value = 6000
$ws.Range("O293").Value = 6500
$ws.Range("P293").Value = 6250
$ws.Range("S293").Value = 625

# Row 294 (was Clemenuless/44399 -> becomes Murcott/44167, Tercera)
$ws.Range("D294").Value = 44167
$ws.Range("K294").Value = "Murcott"
$ws.Range("M294").Value = 200
$ws.Range("N294").Value = 5000
$ws.Range("O294").Value = 5500
$ws.Range("P294").Value = 5250
$ws.Range("S294").Value = 525

# Row 295 (Clemenuless / Especial) date 44400 -> 44399, M 400 -> 500
$ws.Range("D295").Value = 44399
$ws.Range("M295").Value = 500

# Row 296 (Clemenuless / Primera) date 44400 -> 44399, M 500 -> 600
$ws.Range("D296").Value = 44399
$ws.Range("M296").Value = 600

# Row 297 (Clemenuless / Segunda) date 44400 -> 44399 (no other change)
$ws.Range("D297").Value = 44399

# Row 298 (Clemenuless / Tercera) date 44400 -> 44399, M 300 -> 240
$ws.Range("D298").Value = 44399
$ws.Range("M298").Value = 240

# --- Append new rows 299-302, a new week's data (copy of what row 295-298 held before the edit) ---

$ws.Range("A299").Value = 8
$ws.Range("B299").Value = "Terminal La Palmera de La Serena"
$ws.Range("C299").Value = "Coquimbo"
$ws.Range("D299").Value = 44400
$ws.Range("D299").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E299").Value = 4
$ws.Range("F299").Value = "Fruta"
$ws.Range("G299").Value = 100102
$ws.Range("H299").Value = "Cítricos"
$ws.Range("I299").Value = 100102004
$ws.Range("J299").Value = "Mandarina"
$ws.Range("K299").Value = "Clemenuless"
$ws.Range("L299").Value = "Especial"
$ws.Range("M299").Value = 400
$ws.Range("N299").Value = 5500
$ws.Range("O299").Value = 6000
$ws.Range("P299").Value = 5750
$ws.Range("Q299").Value = "$/bandeja 10 kilos"
$ws.Range("R299").Value = "Provincia de Limarí"
$ws.Range("S299").Value = 575
$ws.Range("T299").Value = 10

$ws.Range("A300").Value = 8
$ws.Range("B300").Value = "Terminal La Palmera de La Serena"
$ws.Range("C300").Value = "Coquimbo"
$ws.Range("D300").Value = 44400
$ws.Range("D300").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E300").Value = 4
$ws.Range("F300").Value = "Fruta"
$ws.Range("G300").Value = 100102
$ws.Range("H300").Value = "Cítricos"
$ws.Range("I300").Value = 100102004
$ws.Range("J300").Value = "Mandarina"
$ws.Range("K300").Value = "Clemenuless"
$ws.Range("L300").Value = "Primera"
$ws.Range("M300").Value = 500
$ws.Range("N300").Value = 4500
$ws.Range("O300").Value = 5000
$ws.Range("P300").Value = 4750
$ws.Range("Q300").Value = "$/bandeja 10 kilos"
$ws.Range("R300").Value = "Provincia de Limarí"
$ws.Range("S300").Value = 475
$ws.Range("T300").Value = 10

$ws.Range("A301").Value = 8
$ws.Range("B301").Value = "Terminal La Palmera de La Serena"
$ws.Range("C301").Value = "Coquimbo"
$ws.Range("D301").Value = 44400
$ws.Range("D301").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E301").Value = 4
$ws.Range("F301").Value = "Fruta"
$ws.Range("G301").Value = 100102
$ws.Range("H301").Value = "Cítricos"
$ws.Range("I301").Value = 100102004
$ws.Range("J301").Value = "Mandarina"
$ws.Range("K301").Value = "Clemenuless"
$ws.Range("L301").Value = "Segunda"
$ws.Range("M301").Value = 440
$ws.Range("N301").Value = 3500
$ws.Range("O301").Value = 4000
$ws.Range("P301").Value = 3750
$ws.Range("Q301").Value = "$/bandeja 10 kilos"
$ws.Range("R301").Value = "Provincia de Limarí"
$ws.Range("S301").Value = 375
$ws.Range("T301").Value = 10

$ws.Range("A302").Value = 8
$ws.Range("B302").Value = "Terminal La Palmera de La Serena"
$ws.Range("C302").Value = "Coquimbo"
$ws.Range("D302").Value = 44400
$ws.Range("D302").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E302").Value = 4
$ws.Range("F302").Value = "Fruta"
$ws.Range("G302").Value = 100102
$ws.Range("H302").Value = "Cítricos"
$ws.Range("I302").Value = 100102004
$ws.Range("J302").Value = "Mandarina"
$ws.Range("K302").Value = "Clemenuless"
$ws.Range("L302").Value = "Tercera"
$ws.Range("M302").Value = 300
$ws.Range("N302").Value = 2500
$ws.Range("O302").Value = 3000
$ws.Range("P302").Value = 2750
$ws.Range("Q302").Value = "$/bandeja 10 kilos"
$ws.Range("R302").Value = "Provincia de Limarí"
$ws.Range("S302").Value = 275
$ws.Range("T302").Value = 10
